$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new test results for row 6 (split0005)
$ws.Range("B6").Value = 13335840
$ws.Range("E6").Value = "ok"

# Add new test results for row 12 (split0011)
$ws.Range("B12").Value = 13258080
$ws.Range("E12").Value = "ok"

# Update the active selection to E12 as recorded after the edit
$ws.Range("E12").Select()
